$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "youlchikk"
$ws.Range("B2").Value = "3863055"
$ws.Range("C2").Value = "yulia@mail.ru"
$ws.Range("D2").Value = "06-05-2004"
$ws.Range("E2").Value = "password"

$ws.Columns.AutoFit() | Out-Null

$ws.Range("C11").Select()
